$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# 1. Update every cell's column width across all 9 rows (5 columns).
#    dxa (twips) -> points conversion: points = dxa / 20
# ---------------------------------------------------------------------
$widthsDxa = @(2099, 8060, 1266, 1306, 1217)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Width = $widthsDxa[$c - 1] / 20
    }
}

# ---------------------------------------------------------------------
# 2. Row 4 ("Product Backlog.") - merge the split date runs
#    "2024/1" + "1" + "/" + "22" -> single run "2024/11/22"
# ---------------------------------------------------------------------
$dateCell = $t.Cell(4, 4)
$dateCell.Range.Find.Execute("2024/11/22", $true, $false, $false, $false, $false, $true, 1, $false, "2024/11/22", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Rows 5-7 - add a task-name paragraph before the existing blank
#    paragraph in column 1.
# ---------------------------------------------------------------------
$newTaskNames = @{
    5 = "Use case diagram created."
    6 = "Full use cases defined."
    7 = "Tasks defined as user stories."
}
foreach ($rowNum in $newTaskNames.Keys) {
    $cell = $t.Cell($rowNum, 1)
    $insertRng = $cell.Range.Duplicate
    $insertRng.Collapse(1)
    $insertRng.InsertBefore($newTaskNames[$rowNum] + "`r")
}

# ---------------------------------------------------------------------
# 4. Row 8 ("GitHub project for coursework set-up.") - populate every
#    column with the new content.
# ---------------------------------------------------------------------
# Column 1: Tasks - new paragraph inserted before the existing blank one.
$c1 = $t.Cell(8, 1)
$c1Rng = $c1.Range.Duplicate
$c1Rng.Collapse(1)
$c1Rng.InsertBefore("GitHub project for coursework set-up.`r")

# Column 2: Description - replace the blank paragraph's text.
$c2 = $t.Cell(8, 2)
$c2.Range.Text = "This would include creating the repository along side all the documentation, application files, setting up the dependencies, the Git ignore, creating the Workflow, Learning to Utilise Docker, Git hub issues, Actions, learning to use maven and self-contained Jars, creating the first release branch, "

# Column 4: Due Date - replace the blank paragraph with two paragraphs.
$c4 = $t.Cell(8, 4)
$c4.Range.Text = "2024/11/22`rFriday"

# Column 5: Status - replace the blank paragraph's text.
$c5 = $t.Cell(8, 5)
$c5.Range.Text = "Pending."

# ---------------------------------------------------------------------
# 5. Remove five of the eleven blank paragraphs that sit between the
#    "Software Engineering Methods Project Backlog." heading and the
#    "How to utilise the Backlog." heading.
#
#    NOTE: after touching the table via $d.Tables.Item(...), the
#    $d.Paragraphs collection becomes unreliable (it keeps returning
#    the first paragraph's text for every index). Using
#    $d.Content.Paragraphs instead sidesteps that issue.
# ---------------------------------------------------------------------
$bodyRange = $d.Content
$headingIdx = 0
for ($i = 1; $i -le $bodyRange.Paragraphs.Count; $i++) {
    $ptxt = $bodyRange.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptxt -eq "Software Engineering Methods Project Backlog.") {
        $headingIdx = $i
        break
    }
}
for ($k = 1; $k -le 5; $k++) {
    $p = $d.Content.Paragraphs.Item($headingIdx + 1)
    $p.Range.Delete()
}
